# Weekly update: insert a new price record as row 63, pushing all
# subsequent "Cebollín baby" rows down by one (63->64, ..., 151->152).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 63; Excel shifts rows
# 63..151 down to 64..152 and carries formatting down from row 62.
$ws.Rows.Item(63).Insert()

# Populate the new row with the latest week's record.
$ws.Range("A63").Value = 1
$ws.Range("B63").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C63").Value = 'Arica y Parinacota'
$ws.Range("D63").Value = 45272
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = 100112038
$ws.Range("G63").Value = 'Cebollín baby'
$ws.Range("H63").Value = 'Sin especificar'
$ws.Range("I63").Value = 'Primera'
$ws.Range("J63").Value = 350
$ws.Range("K63").Value = 1800
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = 1886
$ws.Range("N63").Value = '$/paquete 1,5 a 2 kilos'
$ws.Range("O63").Value = 'Región de Arica y Parinacota'
$ws.Range("P63").Value = 943
$ws.Range("Q63").Value = 2
$ws.Range("R63").Value = 'Hortaliza'
